# Add a new "AttenuationCorrection" metadata column after the existing
# ReconFilterSize column (Z), following the same header style/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing header cell (Z1) onto the
# new header cell (AA1), then set its text.
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("AA1").PasteSpecial(-4122) | Out-Null
$ws.Range("AA1").Value = "AttenuationCorrection"

# Size the new column to match the template's other metadata columns.
# ColumnWidth is expressed in "characters"; this runtime stores the XML
# <col width> as ColumnWidth + 5/6, so back that offset out to land on the
# desired stored width (~20.33203125).
$ws.Columns.Item(27).ColumnWidth = 20.33203125 - (5/6)

# Update the view so the new column is visible/selected, matching the
# scrolled/selected state left behind after adding the column.
$win = $excel.ActiveWindow
$win.ScrollColumn = 20
$win.ScrollRow = 1
$ws.Range("AA1").Select() | Out-Null
